# Daily attendance processing - 2025-12-29 06:19:18
#
# Source data refresh for "Y4_B2526_General_&_Special_Surgery_1_B1" session
# analysis:
#   - Group B1E1 gained one student (25 -> 26), so every B1E1 session's
#     "Students" fraction (column H, rows 212-238) and the Group Statistics
#     row for B1E1 (row 23: Students / Avg Attendance %) are recomputed, and
#     the workbook-wide Total Students count (L4) increments accordingly.
#   - The "Recorded By" text (column G) is re-rendered with the contributor
#     list in a new order across many already-recorded sessions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# --- Total Students (Class Statistics) --------------------------------
$ws.Range("L4").Value = 321

# --- Group Statistics row for B1E1 (row 23) ----------------------------
$ws.Range("M23").Value = 26
# "80.3%" reads as text, like the rest of this % column; a leading
# apostrophe keeps Excel from coercing it into a numeric percentage.
$ws.Range("S23").Value = "'80.3%"

# --- Per-session "Students" fraction for group B1E1 (rows 212-238) -----
$b1e1Students = [ordered]@{
    H212 = "20/26"
    H213 = "19/26"
    H214 = "3/26"
    H215 = "24/26"
    H216 = "20/26"
    H217 = "24/26"
    H218 = "26/26"
    H219 = "24/26"
    H220 = "25/26"
    H221 = "25/26"
    H222 = "21/26"
    H223 = "22/26"
    H224 = "22/26"
    H225 = "20/26"
    H226 = "18/26"
    H227 = "0/26"
    H228 = "0/26"
    H229 = "0/26"
    H230 = "0/26"
    H231 = "0/26"
    H232 = "0/26"
    H233 = "0/26"
    H234 = "0/26"
    H235 = "0/26"
    H236 = "0/26"
    H237 = "0/26"
    H238 = "0/26"
}
foreach ($addr in $b1e1Students.Keys) {
    $ws.Range($addr).Value = $b1e1Students[$addr]
}

# --- "Recorded By" (column G): reorder to "dnasr281@gmail.com, System" -
$recordedByRows = @(
    8, 9, 10, 12, 14, 15, 17,
    34, 35, 36, 38, 40, 41, 43,
    60, 61, 62, 64, 66, 67, 69,
    86, 87, 88, 90, 92, 93, 95,
    112, 113, 114, 116, 118, 119, 121,
    138, 139, 140, 142, 144, 145, 147,
    164, 167, 170,
    191, 194, 197,
    218, 221, 224,
    245, 248, 251,
    272, 275, 278,
    299, 302, 305
)
foreach ($row in $recordedByRows) {
    $ws.Range("G$row").Value = "dnasr281@gmail.com, System"
}
